$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "71.031.19"
$ws.Range("E2").Value = "  +1.93%  "

# Row 3
$ws.Range("D3").Value = "3.685.31"
$ws.Range("E3").Value = "  +7.60%  "

# Row 4
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.90"
$ws.Range("E5").Value = "  -0.27%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.66"
$ws.Range("E6").Value = "  +0.55%  "

# Row 7
$ws.Range("D7").Value = "3.670.32"
$ws.Range("E7").Value = "  +7.37%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.617"
$ws.Range("E8").Value = "  +3.96%  "

# Row 9
$ws.Range("E9").Value = "  +0.13%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.199"
$ws.Range("E10").Value = "  -0.83%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.88"
$ws.Range("E11").Value = "  +26.40%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.609"
$ws.Range("E12").Value = "  +4.30%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "49.00"
$ws.Range("E13").Value = "  +0.46%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000287"
$ws.Range("E14").Value = "  +1.58%  "

# Row 15
$ws.Range("D15").Value = "4.284.52"
$ws.Range("E15").Value = "  +7.88%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "678.08"
$ws.Range("E16").Value = "  -2.32%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "9.00"
$ws.Range("E17").Value = "  +4.44%  "

# Row 18
$ws.Range("D18").Value = "3.695.89"
$ws.Range("E18").Value = "  +7.93%  "

# Row 19
$ws.Range("D19").Value = "71.213.14"
$ws.Range("E19").Value = "  +2.22%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.122"
$ws.Range("E20").Value = "  +0.78%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.94"
$ws.Range("E21").Value = "  +1.53%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.57"
$ws.Range("E22").Value = "  +2.02%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.940"
$ws.Range("E23").Value = "  +4.85%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.34"
$ws.Range("E24").Value = "  +2.57%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "101.80"
$ws.Range("E25").Value = "  +0.63%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.97"
$ws.Range("E26").Value = "  +1.54%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.84"
$ws.Range("E27").Value = "  +6.72%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.21"
$ws.Range("E28").Value = "  +6.39%  "

# Row 29
$ws.Range("E29").Value = "  -0.03%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.11"
$ws.Range("E30").Value = "  +4.93%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.42"
$ws.Range("E31").Value = "  +4.84%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.12"
$ws.Range("E32").Value = "  +4.11%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.42"
$ws.Range("E33").Value = "  -2.37%  "

# Row 34
$ws.Range("E34").Value = "  +5.25%  "

# Row 35
$ws.Range("E35").Value = "  +8.62%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "580.89"
$ws.Range("E36").Value = "  +1.50%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.18"
$ws.Range("E37").Value = "  +1.69%  "

# Row 38
$ws.Range("E38").Value = "  +4.20%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "58.68"
$ws.Range("E39").Value = "  +0.79%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.03%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0462"
$ws.Range("E41").Value = "  +10.53%  "

# Row 42
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.620.56"
$ws.Range("E42").Value = "  +1.69%  "

# Row 43
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.144"
$ws.Range("E43").Value = "  +3.35%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.351"
$ws.Range("E44").Value = "  +5.56%  "

# Row 45
$ws.Range("D45").Value = "0.0₃0765"
$ws.Range("E45").Value = "  +4.16%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "35.24"
$ws.Range("E46").Value = "  +0.35%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.76"
$ws.Range("E47").Value = "  +3.20%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.89"
$ws.Range("E48").Value = "  +9.51%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.133"
$ws.Range("E49").Value = "  +3.99%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.59"
$ws.Range("E50").Value = "  +1.28%  "

# Row 51
$ws.Range("E51").Value = "  +10.30%  "
